# Attendance correction (12 March update): the "A" (absent) mark in the
# session-19 column (Y) is removed for a batch of students. Clearing that
# mark also drops the "Total Absence" count in column E by one, since E
# is =IF(D>0,COUNTIF(G:Z,"A"),"") and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(11, 17, 18, 20, 22, 25, 29, 31, 35, 37, 40, 42, 44, 45, 46, 47, 48, 50, 51, 54, 55, 58, 60, 61, 62, 63, 64, 65, 67, 70, 74, 75, 76, 78, 79)

foreach ($r in $rows) {
    # Drop the "A" value itself...
    $ws.Range("Y$r").ClearContents()
    # ...then restore the plain (unmarked) cell formatting by pulling it
    # from the neighbouring already-blank cell in the same row, so the
    # cell's style matches every other untouched attendance cell.
    $ws.Range("Z$r").Copy() | Out-Null
    $ws.Range("Y$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
$excel.CalculateFullRebuild()
